$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 35: add "Saving project is quite messed up." in column B (alongside existing A35 text)
$ws.Range("B35").Value = "Saving project is quite messed up."

# Row 36: add "Saving project is quite messed up." in column B (alongside existing A36 text)
$ws.Range("B36").Value = "Saving project is quite messed up."

# Row 37 (old A37 held the same "Saving project is quite messed up." text) is no longer
# needed now that its content lives in B35/B36 - remove the row entirely.
$ws.Range("A37:B37").Delete()

# Update the selected cell to match the new edit location.
$ws.Range("B28").Select()
